# Rename the (only) worksheet from the generic default "Лист1" to a
# descriptive name matching the workbook's file name, and move the
# selection to where the author last left it (B20) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "stepen_vnedreniya_ikt"

# Restore the cell selection/cursor position that was saved with the
# workbook (was Q10, now B20).
$ws.Range("B20").Select()
